$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '50.961.93'
$ws.Range('E2').Value = '  -1.07%  '

# Row 3
$ws.Range('D3').Value = '2.933.26'
$ws.Range('E3').Value = '  -1.80%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.13%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '374.10'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.67%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.05%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.535'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.27%  '

# Row 8
$ws.Range('E8').Value = '  +0.21%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.581'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.66%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.88'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.25%  '

# Row 11
$ws.Range('E11').Value = '  -0.70%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0843'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.18%  '

# Row 13
$ws.Range('D13').Value = '3.394.62'
$ws.Range('E13').Value = '  -1.23%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.24'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +65.00%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '17.87'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.39%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.38'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.73%  '

# Row 17
$ws.Range('D17').Value = '2.941.34'
$ws.Range('E17').Value = '  -1.04%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.968'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.29%  '

# Row 19
$ws.Range('D19').Value = '50.893.24'
$ws.Range('E19').Value = '  -1.07%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.68%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.38'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.61%  '

# Row 22
$ws.Range('D22').Value = '0.0₃0952'
$ws.Range('E22').Value = '  -1.26%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '264.68'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.97%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.72%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.17'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +12.14%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.44%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.48'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.36%  '

# Row 28
$ws.Range('E28').Value = '  +0.08%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.166'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.96%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '25.46'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.01%  '

# Row 31
$ws.Range('E31').Value = '  -3.45%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.95'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.27%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '50.53'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.58%  '

# Row 34
$ws.Range('E34').Value = '  -2.99%  '

# Row 35
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '32.77'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.42%  '

# Row 36
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0438'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.22%  '

# Row 37
$ws.Range('E37').Value = '  -0.16%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.17'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.16%  '

# Row 39
$ws.Range('E39').Value = '  -0.75%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.16'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.60%  '

# Row 41
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.77'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.29%  '

# Row 42
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.46'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.73%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '120.26'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.91%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.95'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.74%  '

# Row 45
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.274'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.95%  '

# Row 46
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.95%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.28'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.70%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.30'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.91%  '

# Row 49
$ws.Range('D49').Value = '2.002.30'
$ws.Range('E49').Value = '  -2.30%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0329'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.51%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.29'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.18%  '
